$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$cell = $t.Cell(4, 3)
$rng = $cell.Range

# Cell shading: clear pattern, auto foreground, D86DCB background (accent5, lighter 60%)
$cell.Shading.Texture = 0
$cell.Shading.ForegroundPatternColor = -16777216
$cell.Shading.BackgroundPatternColor = 13331928

# Run + paragraph-mark formatting: bold, 24pt (48 half-points), accent3-lighter-60% green
$rng.Font.Bold = 1
$rng.Font.BoldBi = 1
$rng.Font.Size = 24
$rng.Font.SizeBi = 24
$rng.Font.Color = 5887047

Write-Output "done"
